$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New label / id cells -------------------------------------------------
# Order chosen to reproduce the author's shared-string insertion order.
$ws.Range("B1").Value = "id=priceTable"

$ws.Range("A5").Value = "Choose Silver"
$ws.Range("A6").Value = "Choose Gold"
$ws.Range("A7").Value = "Choose Platinum"
$ws.Range("A8").Value = "Choose Ultimate"

$ws.Range("J5").Value = "X"
$ws.Range("J6").Value = "X"
$ws.Range("J7").Value = "X"
$ws.Range("J8").Value = "X"

$ws.Range("C1").Value = "id=selectsilver"
$ws.Range("D1").Value = "id=selectgold"
$ws.Range("F1").Value = "id=selectultimate"
$ws.Range("G1").Value = "id=viewquote"
$ws.Range("H1").Value = "id=downloadquote"
$ws.Range("I1").Value = "id=preventerproductdata"
$ws.Range("J1").Value = "id=nextsendquote"

$ws.Range("E1").Value = "*css=label >> css=[value=Platinum]"

$ws.Range("C5").Value = "<CHECK>"
$ws.Range("D6").Value = "<CHECK>"
$ws.Range("E7").Value = "<CHECK>"
$ws.Range("F8").Value = "<CHECK>"

# --- Yellow highlight fill on the spacer rows -----------------------------
$ws.Range("B2:J2").Interior.Color = 65535
$ws.Range("B4:J4").Interior.Color = 65535

# --- Column width changes ---------------------------------------------------
# (ColumnWidth values below are chosen so the engine's internal pixel
# quantization lands on the closest achievable width to the authored target.)
$ws.Columns.Item(5).ColumnWidth = 36.8640625
$ws.Columns.Item(9).ColumnWidth = 28.8640625
$ws.Columns.Item(10).ColumnWidth = 24.00671875

# --- Reposition / resize the picture ---------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 135.8351181102362
$shp.Left = 0
$shp.Width = 865.8
$shp.Height = 501.2827559055118

# --- Selection shown when the file was last saved ---------------------------
$ws.Range("I12").Select()
